# Update "想去人数" (F column) values in the "展览" and "全部类型" sheets
# to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 3238
    7  = 1703
    8  = 1639
    17 = 233
    23 = 384
    29 = 319
    30 = 2217
    34 = 442
    37 = 229
    38 = 348
    40 = 524
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
